$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the last data row (132) to make room for a
# new time-log entry plus a blank separator row before the summary rows.
$ws.Rows("133:134").Insert() | Out-Null

# Row 132: this was previously a mis-typed "2ß14" text value in A132 with no
# end-time, fix the year to a real number and supply the missing end time.
$ws.Range("A132").Value = 2014
$ws.Range("E132").Value = 0.92708333333333337

# New row 133: another entry for the same day.
$ws.Range("A133").Value = 2014
$ws.Range("B133").Value = 7
$ws.Range("C133").Value = 12
$ws.Range("D133").Value = 0.375
$ws.Range("E133").Value = 0.5
$ws.Range("F133").Formula = "=(E133-D133)*24*60"
$ws.Range("G133").Formula = "=F133/60"

# The row insert copied G132's format down into the new blank row 134 - that
# row should stay fully empty (only D/E/F keep the blank formatted cells).
$ws.Range("G134").Clear() | Out-Null

# Re-point the summary formulas at the new last data row.
$ws.Range("F135").Formula = "=SUM(F2:F133)"

# Leave the same cell selected as in the authored workbook.
$ws.Range("F133").Select() | Out-Null
